$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(11, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/12/2025"
$dateCell.ClearFormats()

$ws.Cells.Item(11, 2).Value = 0.1205199642338759
$ws.Cells.Item(11, 3).Value = 0.8794800357661241
